$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "DP" (sheet2.xml): add a "餐點"/"葷素" classification column, insert
# the missing meal rows (so it mirrors the full 12-item list from "DP (2)"),
# and rename two vegetarian-disguised items.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("DP")

# Make room for five new rows (they will become rows 6-10); existing rows
# 6-7 (魷魚羹/功夫麵) shift down to rows 11-12.
$ws2.Rows("6:10").Insert()

# Header row
$ws2.Range("A1").Value = "餐點"
$ws2.Range("F1").Value = "葷素"

# Row 2 - 叻沙牛肉麵 (unchanged) -> meat
$ws2.Range("F2").Value = "葷"

# Row 3 - was 石鍋拌飯, now disguised vegetarian item -> veg
$ws2.Range("A3").Value = "我是素食"
$ws2.Range("F3").Value = "素"

# Row 4 - was 花生雞腿堡, now disguised vegetarian item -> veg
$ws2.Range("A4").Value = "神奇沙拉"
$ws2.Range("F4").Value = "素"

# Row 5 - 原汁牛肉麵 (unchanged) -> meat
$ws2.Range("F5").Value = "葷"

# Row 6 (new) - 炸豬排飯
$ws2.Range("A6").Value = "炸豬排飯"
$ws2.Range("B6").Value = 570
$ws2.Range("C6").Value = 110
$ws2.Range("D6").Value = 71
$ws2.Range("E6").Value = 67
$ws2.Range("F6").Value = "葷"

# Row 7 (new) - 咖哩飯
$ws2.Range("A7").Value = "咖哩飯"
$ws2.Range("B7").Value = 550
$ws2.Range("C7").Value = 120
$ws2.Range("D7").Value = 69
$ws2.Range("E7").Value = 65
$ws2.Range("F7").Value = "葷"

# Row 8 (new) - 雞胸肉餐
$ws2.Range("A8").Value = "雞胸肉餐"
$ws2.Range("B8").Value = 450
$ws2.Range("C8").Value = 135
$ws2.Range("D8").Value = 68
$ws2.Range("E8").Value = 66
$ws2.Range("F8").Value = "葷"

# Row 9 (new) - 素肉全餐
$ws2.Range("A9").Value = "素肉全餐"
$ws2.Range("B9").Value = 440
$ws2.Range("C9").Value = 105
$ws2.Range("D9").Value = 65
$ws2.Range("E9").Value = 63
$ws2.Range("F9").Value = "素"

# Row 10 (new) - 涼麵(素)
$ws2.Range("A10").Value = "涼麵(素)"
$ws2.Range("B10").Value = 380
$ws2.Range("C10").Value = 60
$ws2.Range("D10").Value = 62
$ws2.Range("E10").Value = 61
$ws2.Range("F10").Value = "素"

# Row 11 (was row 6) - 魷魚羹 -> meat
$ws2.Range("F11").Value = "葷"

# Row 12 (was row 7) - 功夫麵 -> meat
$ws2.Range("F12").Value = "葷"

# Row 13 (new) - 烤雞排
$ws2.Range("A13").Value = "烤雞排"
$ws2.Range("B13").Value = 300
$ws2.Range("C13").Value = 60
$ws2.Range("D13").Value = 50
$ws2.Range("E13").Value = 50
$ws2.Range("F13").Value = "葷"

# Column widths: A and E become custom widths, matching the widened table.
$ws2.Columns.Item(1).ColumnWidth = 15
$ws2.Columns.Item(5).ColumnWidth = 13.285714285714286

# ---------------------------------------------------------------------------
# Sheet "DP (2)" (sheet1.xml): view-state selection only moved from G13 to
# E13, still covering A1:E13.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("DP (2)")
$ws1.Range("A1:E13").Select() | Out-Null

# ---------------------------------------------------------------------------
# Restore "DP" as the active sheet/tab with its own selection, since it was
# the active tab before and after the edit.
# ---------------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("J17").Select() | Out-Null
